# GDP.xlsx data_import rework:
# - drop the "Total" column (S), which held a row-wise SUM(B:R) check
# - drop the derived percentage-of-total row (63), which divided each
#   country's row-62 figure by the row-62 total
# - leave the selection on Q63 (last cell of the now-removed share row)
#   as that's where the cursor ended up after the edit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column S (19th column) - the "Total" SUM(B:R) column
$ws.Columns.Item(19).Delete()

# Row 63 - the row of "value / total" percentage formulas
$ws.Rows.Item(63).Delete()

# Match the author's final selection
$ws.Range("Q63").Select()

# Re-saved from an English-locale Excel, so the builtin cell style's
# display name comes out as "Normal" instead of the Dutch "Standaard"
$wb.Styles.Item("Standaard").Delete()
$wb.Styles.Add("Normal")
